$d = $word.ActiveDocument

# Locate the target sentence dynamically so we don't depend on a hardcoded
# paragraph index.
$target = "How do your findings in questions 2 and 3 relate?"
$hit = $d.Content
$found = $hit.Find.Execute($target, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate target sentence"
}

$base = $hit.Start

# Character offsets (relative to $base) of the two digits that change:
#   "How do your findings in questions " [0,34)
#   "2"                                   [34,35)  -> becomes "3"
#   " and "                               [35,40)
#   "3"                                   [40,41)  -> becomes "4"
#   " relate?"                            [41,49)

# Step 1: rewrite the digits in place while the sentence is still a single
# run (mutating .Text on a sub-range does not split the run).
$num1 = $d.Range($base + 34, $base + 35)
$num1.Text = "3"

$num2 = $d.Range($base + 40, $base + 41)
$num2.Text = "4"

# Step 2: force the run to split into five separate runs at the digit
# boundaries by nudging (and reverting) direct character formatting on
# each narrow sub-range. Touching formatting on a sub-range causes the
# host to break the run apart while leaving the run properties unchanged.
$num1 = $d.Range($base + 34, $base + 35)
$num1.Font.Bold = $true
$num1.Font.Bold = $false

$num2 = $d.Range($base + 40, $base + 41)
$num2.Font.Bold = $true
$num2.Font.Bold = $false

$mid = $d.Range($base + 35, $base + 40)
$mid.Font.Bold = $true
$mid.Font.Bold = $false

$tail = $d.Range($base + 41, $base + 49)
$tail.Font.Bold = $true
$tail.Font.Bold = $false

Write-Host "Final text: $($d.Range($base, $base + 49).Text)"
